# Apply cryptos list update (commit: "Updated cryptos list on Fri Sep 29 16:46:15 UTC 2023 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.857.25"
$ws.Range("E2").Value = "  -0.90%  "
$ws.Range("D3").Value = "1.667.99"
$ws.Range("E3").Value = "  +0.79%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").Value = "'215.48"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").Value = "'0.529"
$ws.Range("E6").Value = "  +4.04%  "
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("E8").Value = "  +1.68%  "
$ws.Range("E9").Value = "  +0.74%  "
$ws.Range("D10").Value = "'20.19"
$ws.Range("E10").Value = "  +3.48%  "
$ws.Range("D11").Value = "'0.0895"
$ws.Range("E11").Value = "  +3.97%  "
$ws.Range("D12").Value = "1.903.86"
$ws.Range("E12").Value = "  +0.72%  "
$ws.Range("D13").Value = "1.707.81"
$ws.Range("E13").Value = "  +3.15%  "
$ws.Range("E14").Value = "  +0.30%  "
$ws.Range("E15").Value = "  +1.02%  "
$ws.Range("D16").Value = "'65.84"
$ws.Range("E16").Value = "  +1.40%  "
$ws.Range("D17").Value = "26.883.89"
$ws.Range("E17").Value = "  -0.78%  "
$ws.Range("D18").Value = "'231.62"
$ws.Range("E18").Value = "  -3.78%  "
$ws.Range("E19").Value = "  -0.94%  "
$ws.Range("E20").Value = "  +0.64%  "
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("D22").Value = "'4.46"
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("D23").Value = "'2.21"
$ws.Range("E23").Value = "  -2.09%  "
$ws.Range("D24").Value = "'9.19"
$ws.Range("E24").Value = "  -0.69%  "
$ws.Range("D25").Value = "'145.59"
$ws.Range("E25").Value = "  -0.29%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "'7.13"
$ws.Range("E26").Value = "  +0.22%  "
$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").Value = "'0.116"
$ws.Range("E27").Value = "  +1.34%  "
$ws.Range("D28").Value = "'15.89"
$ws.Range("E28").Value = "  +0.28%  "
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").Value = "'0.0496"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("E31").Value = "  +0.64%  "
$ws.Range("E32").Value = "  +1.55%  "
$ws.Range("D33").Value = "1.464.51"
$ws.Range("E33").Value = "  -3.71%  "
$ws.Range("D34").Value = "'3.15"
$ws.Range("E34").Value = "  +3.60%  "
$ws.Range("E35").Value = "  +3.91%  "
$ws.Range("E36").Value = "  -0.19%  "
$ws.Range("D37").Value = "'0.897"
$ws.Range("E37").Value = "  +0.79%  "
$ws.Range("D38").Value = "'0.570"
$ws.Range("E38").Value = "  -1.34%  "
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("D40").Value = "'5.80"
$ws.Range("E40").Value = "  -2.64%  "
$ws.Range("E41").Value = "  +0.35%  "
$ws.Range("E42").Value = "  +0.92%  "
$ws.Range("D43").Value = "'0.976"
$ws.Range("E43").Value = "  +7.07%  "
$ws.Range("D44").Value = "'65.66"
$ws.Range("E44").Value = "  +1.49%  "
$ws.Range("D45").Value = "1.813.34"
$ws.Range("E45").Value = "  +0.91%  "
$ws.Range("E46").Value = "  +1.29%  "
$ws.Range("D47").Value = "'90.21"
$ws.Range("E47").Value = "  -0.35%  "
$ws.Range("E48").Value = "  -0.54%  "
$ws.Range("E49").Value = "  +2.43%  "
$ws.Range("E50").Value = "  +0.97%  "
$ws.Range("E51").Value = "  +0.66%  "
